$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 255, pushing the existing rows 255-342 down to 256-343.
$ws.Rows(255).Insert()

# Populate the newly inserted row 255 with the new record.
$ws.Range("A255").Value = 10
$ws.Range("B255").Value = "Vega Modelo de Temuco"
$ws.Range("C255").Value = "La Araucanía"
$ws.Range("D255").Value = 45215
$ws.Range("E255").Value = 9
$ws.Range("F255").Value = 100112005
$ws.Range("G255").Value = "Puerro"
$ws.Range("H255").Value = "Azul de Maquehue"
$ws.Range("I255").Value = "Primera"
$ws.Range("J255").Value = 110
$ws.Range("K255").Value = 8000
$ws.Range("L255").Value = 8000
$ws.Range("M255").Value = 8000
$ws.Range("N255").Value = "$/docena de paquetes"
$ws.Range("O255").Value = "Provincia de Cautín"
$ws.Range("P255").Value = 667
$ws.Range("Q255").Value = 12
$ws.Range("R255").Value = "Hortaliza"
